# Working version of eroski with Retrievers iterator class
# Rename the sheet to reflect the new pricing date, then update the
# Eroski (C) / BM (D) price columns for the products that changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet from the old date to the new date.
$ws.Name = "2022-11-19"

function Set-Price($Row, $Value) {
    $ws.Cells.Item($Row, 3).Value = $Value
    $ws.Cells.Item($Row, 4).Value = $Value
}

# Price updates (Eroski/BM columns kept in sync, as in the source data).
Set-Price 8   1.96    # AJO
Set-Price 17  3.05    # AVELLANA TOSTADA PELADA SIN SA
Set-Price 32  13.45   # CERVEZA
Set-Price 40  1.99    # DENTRIFICO
Set-Price 62  1.99    # JUDIA VERDE PLANA
Set-Price 65  6.99    # LANGOSTINO
Set-Price 68  1.09    # LECHE DESNATADA
Set-Price 69  1.09    # LECHE SEMIDESNATADA (was blank)
Set-Price 70  1.19    # LECHUGA
Set-Price 75  3.99    # LIMPIADOR SUELO
Set-Price 90  16.43   # MERLUZA FRESCA ENTERA
Set-Price 92  1.69    # NARANJA
Set-Price 93  1.75    # NATA COCINA (was blank)
Set-Price 94  3.99    # NUEZ
Set-Price 108 2       # PIMIENTO ROJO
Set-Price 109 3.98    # PIÑA
Set-Price 122 4.27    # SALSA BECHAMEL (was blank)
Set-Price 127 1.28    # TOMATE ENSALADA
Set-Price 128 1.7     # TOMATE FRITO ACEITE OLIVA RECETA ARTESANA

# Cells that become blank (cleared, no price available this round).
Set-Price 13  ""      # ANCHOA FILETE ACEITE OLIVA
Set-Price 20  ""      # BONITO
Set-Price 49  ""      # FAIRY
Set-Price 87  ""      # MELOCOTON
Set-Price 100 ""      # PATATA ENTERA CONSERVA (was 1.35, now blank)
